$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -10.15409240999935

$ws.Range("J25").Value = -7.349037835258713
$ws.Range("K25").Value = -1.996536364888143

$ws.Range("I26").Value = -7.759011788180499
$ws.Range("J26").Value = -2.406510317809929
$ws.Range("K26").Value = -2.512938874122935

$ws.Range("H27").Value = -8.324416066900824
$ws.Range("I27").Value = -2.971914596530254
$ws.Range("J27").Value = -3.07834315284326
$ws.Range("K27").Value = -1.329627155229175

$ws.Range("G28").Value = -8.359011788180499
$ws.Range("H28").Value = -3.006510317809929
$ws.Range("I28").Value = -3.112938874122935
$ws.Range("J28").Value = -1.36422287650885
$ws.Range("K28").Value = 1.261604321912398

$ws.Range("F29").Value = -9.333026140304781
$ws.Range("G29").Value = -3.980524669934212
$ws.Range("H29").Value = -4.086953226247218
$ws.Range("I29").Value = -2.338237228633133
$ws.Range("J29").Value = 0.2875899697881152
$ws.Range("K29").Value = -5.355989324100946

$ws.Range("E30").Value = -8.990342712894432
$ws.Range("F30").Value = -3.637841242523862
$ws.Range("G30").Value = -3.744269798836868
$ws.Range("H30").Value = -1.995553801222783
$ws.Range("I30").Value = 0.6302733971984651
$ws.Range("J30").Value = -5.013305896690596
$ws.Range("K30").Value = -0.8595623955320235

$ws.Range("D31").Value = -7.787054418971053
$ws.Range("E31").Value = -2.434552948600483
$ws.Range("F31").Value = -2.540981504913489
$ws.Range("G31").Value = -0.792265507299404
$ws.Range("H31").Value = 1.833561691121844
$ws.Range("I31").Value = -3.810017602767218
$ws.Range("J31").Value = 0.3437258983913554
$ws.Range("K31").Value = -0.3221206215403352

$ws.Range("C32").Value = -10.6590117881805
$ws.Range("D32").Value = -5.306510317809929
$ws.Range("E32").Value = -5.412938874122935
$ws.Range("F32").Value = -3.66422287650885
$ws.Range("G32").Value = -1.038395678087602
$ws.Range("H32").Value = -6.681974971976663
$ws.Range("I32").Value = -2.528231470818091
$ws.Range("J32").Value = -3.194077990749781
$ws.Range("K32").Value = -3.951141586996604

$ws.Range("B33").Value = -25.08470868489006
$ws.Range("C33").Value = -19.73220721451949
$ws.Range("D33").Value = -19.8386357708325
$ws.Range("E33").Value = -18.08991977321841
$ws.Range("F33").Value = -15.46409257479716
$ws.Range("G33").Value = -21.10767186868623
$ws.Range("H33").Value = -16.95392836752765
$ws.Range("I33").Value = -17.61977488745934
$ws.Range("J33").Value = -18.37683848370617
$ws.Range("K33").Value = -14.68761087779509

$ws.Range("B34").Value = 5.352501470370569
$ws.Range("C34").Value = 5.246072914057564
$ws.Range("D34").Value = 6.994788911671648
$ws.Range("E34").Value = 9.620616110092897
$ws.Range("F34").Value = 3.977036816203835
$ws.Range("G34").Value = 8.130780317362408
$ws.Range("H34").Value = 7.464933797430717
$ws.Range("I34").Value = 6.707870201183894
$ws.Range("J34").Value = 10.39709780709497
$ws.Range("K34").Value = 9.963901476152344

$ws.Range("B35").Value = -0.1064285563130061
$ws.Range("C35").Value = 1.642287441301079
$ws.Range("D35").Value = 4.268114639722327
$ws.Range("E35").Value = -1.375464654166735
$ws.Range("F35").Value = 2.778278846991838
$ws.Range("G35").Value = 2.112432327060148
$ws.Range("H35").Value = 1.355368730813325
$ws.Range("I35").Value = 5.044596336724396
$ws.Range("J35").Value = 4.611400005781775
$ws.Range("K35").Value = 3.714584630134098

$ws.Range("B36").Value = 1.748715997614085
$ws.Range("C36").Value = 4.374543196035333
$ws.Range("D36").Value = -1.269036097853729
$ws.Range("E36").Value = 2.884707403304844
$ws.Range("F36").Value = 2.218860883373154
$ws.Range("G36").Value = 1.461797287126331
$ws.Range("H36").Value = 5.151024893037402
$ws.Range("I36").Value = 4.71782856209478
$ws.Range("J36").Value = 3.821013186447104
$ws.Range("K36").Value = 5.08770374598177

$ws.Range("B37").Value = 2.625827198421248
$ws.Range("C37").Value = -3.017752095467813
$ws.Range("D37").Value = 1.135991405690759
$ws.Range("E37").Value = 0.4701448857590689
$ws.Range("F37").Value = -0.286918710487754
$ws.Range("G37").Value = 3.402308895423317
$ws.Range("H37").Value = 2.969112564480695
$ws.Range("I37").Value = 2.072297188833019
$ws.Range("J37").Value = 3.338987748367685
$ws.Range("K37").Value = 3.74600879241737

$ws.Range("B38").Value = -5.643579293889061
$ws.Range("C38").Value = -1.489835792730489
$ws.Range("D38").Value = -2.155682312662179
$ws.Range("E38").Value = -2.912745908909002
$ws.Range("F38").Value = 0.7764816970020689
$ws.Range("G38").Value = 0.3432853660594475
$ws.Range("H38").Value = -0.5535300095882292
$ws.Range("I38").Value = 0.7131605499464369
$ws.Range("J38").Value = 1.120181593996122
$ws.Range("K38").Value = -1.176625030627534

$ws.Range("B39").Value = 4.153743501158573
$ws.Range("C39").Value = 3.487896981226882
$ws.Range("D39").Value = 2.730833384980059
$ws.Range("E39").Value = 6.42006099089113
$ws.Range("F39").Value = 5.986864659948509
$ws.Range("G39").Value = 5.090049284300832
$ws.Range("H39").Value = 6.356739843835499
$ws.Range("I39").Value = 6.763760887885184
$ws.Range("J39").Value = 4.466954263261528
$ws.Range("K39").Value = 5.055752123058397

$ws.Range("B40").Value = -0.6658465199316905
$ws.Range("C40").Value = -1.422910116178513
$ws.Range("D40").Value = 2.266317489732558
$ws.Range("E40").Value = 1.833121158789936
$ws.Range("F40").Value = 0.9363057831422594
$ws.Range("G40").Value = 2.202996342676926
$ws.Range("H40").Value = 2.610017386726611
$ws.Range("I40").Value = 0.3132107621029547
$ws.Range("J40").Value = 0.9020086218998244
$ws.Range("K40").Value = 0.7731291372494979

$ws.Range("B41").Value = -0.7570635962468228
$ws.Range("C41").Value = 2.932164009664248
$ws.Range("D41").Value = 2.498967678721627
$ws.Range("E41").Value = 1.60215230307395
$ws.Range("F41").Value = 2.868842862608616
$ws.Range("G41").Value = 3.275863906658301
$ws.Range("H41").Value = 0.9790572820346453
$ws.Range("I41").Value = 1.567855141831515
$ws.Range("J41").Value = 1.438975657181188
$ws.Range("K41").Value = 1.046227481711355

$ws.Range("B42").Value = 3.689227605911071
$ws.Range("C42").Value = 3.256031274968449
$ws.Range("D42").Value = 2.359215899320773
$ws.Range("E42").Value = 3.625906458855439
$ws.Range("F42").Value = 4.032927502905124
$ws.Range("G42").Value = 1.736120878281468
$ws.Range("H42").Value = 2.324918738078338
$ws.Range("I42").Value = 2.196039253428011
$ws.Range("J42").Value = 1.803291077958178
$ws.Range("K42").Value = 1.276334624211031

$ws.Range("B43").Value = -0.4331963309426214
$ws.Range("C43").Value = -1.330011706590298
$ws.Range("D43").Value = -0.06332114705563208
$ws.Range("E43").Value = 0.343699896994053
$ws.Range("F43").Value = -1.953106727629603
$ws.Range("G43").Value = -1.364308867832733
$ws.Range("H43").Value = -1.49318835248306
$ws.Range("I43").Value = -1.885936527952893
$ws.Range("J43").Value = -2.41289298170004
$ws.Range("K43").Value = -2.98034896160172

$ws.Range("B44").Value = -0.8968153756476767
$ws.Range("C44").Value = 0.3698751838869894
$ws.Range("D44").Value = 0.7768962279366745
$ws.Range("E44").Value = -1.519910396686981
$ws.Range("F44").Value = -0.9311125368901116
$ws.Range("G44").Value = -1.059992021540438
$ws.Range("H44").Value = -1.452740197010272
$ws.Range("I44").Value = -1.979696650757419
$ws.Range("J44").Value = -2.547152630659099

$ws.Range("B45").Value = 1.266690559534666
$ws.Range("C45").Value = 1.673711603584351
$ws.Range("D45").Value = -0.6230950210393047
$ws.Range("E45").Value = -0.03429716124243498
$ws.Range("F45").Value = -0.1631766458927615
$ws.Range("G45").Value = -0.5559248213625949
$ws.Range("H45").Value = -1.082881275109742
$ws.Range("I45").Value = -1.650337255011422

$ws.Range("B46").Value = 0.4070210440496851
$ws.Range("C46").Value = -1.889785580573971
$ws.Range("D46").Value = -1.300987720777101
$ws.Range("E46").Value = -1.429867205427428
$ws.Range("F46").Value = -1.822615380897261
$ws.Range("G46").Value = -2.349571834644408
$ws.Range("H46").Value = -2.917027814546088

$ws.Range("B47").Value = -2.296806624623656
$ws.Range("C47").Value = -1.708008764826786
$ws.Range("D47").Value = -1.836888249477113
$ws.Range("E47").Value = -2.229636424946946
$ws.Range("F47").Value = -2.756592878694093
$ws.Range("G47").Value = -3.324048858595773

$ws.Range("B48").Value = 0.5887978597968697
$ws.Range("C48").Value = 0.4599183751465432
$ws.Range("D48").Value = 0.06717019967670979
$ws.Range("E48").Value = -0.4597862540704372
$ws.Range("F48").Value = -1.027242233972117

$ws.Range("B49").Value = -0.1288794846503265
$ws.Range("C49").Value = -0.5216276601201599
$ws.Range("D49").Value = -1.048584113867307
$ws.Range("E49").Value = -1.616040093768987

$ws.Range("B50").Value = -0.3927481754698334
$ws.Range("C50").Value = -0.9197046292169804
$ws.Range("D50").Value = -1.48716060911866

$ws.Range("B51").Value = -0.5269564537471469
$ws.Range("C51").Value = -1.094412433648827

$ws.Range("B52").Value = -0.5674559799016801
